$d = $word.ActiveDocument
$d.Content.Find.Execute("19 November 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 November 2022", 2)
